$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new value parses as a plain number must be
# forced to Text format first, otherwise Excel auto-converts the
# entry to a numeric value (losing the literal string / trailing zeros).
$ws.Range("D2").Value = "29.201.50"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "1.856.02"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7077"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.61"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07958"
$ws.Range("E8").Value = "  +3.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3023"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.40"
$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08172"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").Value = "1.853.13"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.162"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7014"
$ws.Range("E14").Value = "  -3.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.60"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "29.179.39"
$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.806"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007860"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.22"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.05"
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "2.106.95"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.370"
$ws.Range("E24").Value = "  -3.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.87"
$ws.Range("E25").Value = "  +0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.907"
$ws.Range("E26").Value = "  -1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1430"
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.923"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("E30").Value = "  +2.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.480"
$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.361"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.010"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05194"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("E35").Value = "  -2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6994"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.669"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01838"
$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  +1.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9301"
$ws.Range("E41").Value = "  -0.69%  "

$ws.Range("D42").Value = "1.122.08"
$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4242"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.827"
$ws.Range("E44").Value = "  -3.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.54"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.78"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.755"
$ws.Range("E48").Value = "  -1.38%  "

$ws.Range("D49").Value = "1.989.86"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.097"
$ws.Range("E50").Value = "  -1.19%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.948"
$ws.Range("E51").Value = "  -1.35%  "
